# "Generate Report for Handback" -- refresh the localization-status report
# after a handback run: status flips from "ready for handoff" to "handed
# back", the per-language handback timestamps advance, and the (now
# resolved) stale-handback error details are cleared.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: both language Status columns move from "Ready for
# handoff" to "Handed back: in sync with en-US".
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"

# zh-cn sheet: new handback datetime, error detail resolved (cleared).
$zhcn.Range("K2").Value = "2016-08-17 20:50:32"
$zhcn.Range("P2").Value = ""

# de-de sheet: new handback datetime, error detail resolved (cleared).
$dede.Range("K2").Value = "2016-08-17 20:50:40"
$dede.Range("P2").Value = ""

# Column widths widen/narrow to fit the new text (the Status columns now
# hold longer text; the Error Detail columns are now empty).
$overview.Columns("E:F").ColumnWidth = 29.166666666666668

$zhcn.Columns("C:C").ColumnWidth = 29.166666666666668
$zhcn.Columns("P:P").ColumnWidth = 12.833333333333334

$dede.Columns("C:C").ColumnWidth = 29.166666666666668
$dede.Columns("P:P").ColumnWidth = 12.833333333333334
